# Add materials for session 04
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Woche 6): add "Vorbereitung" (D7) -- entered first so its shared
# string lands at index 35, matching the target shared-strings order.
$ws.Range("D7").Value = "prep/p05.html"

# Row 6 (Woche 5): add "Abgaben" (F6) and "Folien" (E6) for session 04 materials
$ws.Range("F6").Value = "exercises/e04.html"
$ws.Range("E6").Value = "slides/slides.html#/sitzung-05-warum-wir-nicht-alle-medienerlebnisse-als-unterhaltung-bezeichnen-sollten"

# Update the active selection to match the edited cell
$ws.Range("E6").Select()
